# #5: property boat&car done
#
# Rebuild the "汽車" (car) sheet (3rd worksheet) so that:
#  - Row 1 becomes the standard header row used across the other property
#    sheets: name / capacity / owner / register_date / register_reason /
#    acquire_value / property_category / category / date / legislator_name /
#    legislator_id / source_file / index  (columns B..N)
#  - Rows 2-5 (the existing car records) keep their original data in B..G
#    and gain the same trailing metadata columns H..N that every other
#    sheet already has (property_category, category, date, legislator_name,
#    legislator_id, source_file, index).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(3)

# ---- Give the new columns (H..N) the same look as the existing ones ----
# Row 1 uses the bold/bordered header style (copy it from an existing
# header cell); rows 2-5 use the plain data style.
$ws.Cells.Item(1, 2).Copy()
$ws.Range("H1:N1").PasteSpecial(-4122)

$ws.Cells.Item(2, 2).Copy()
$ws.Range("H2:N2").PasteSpecial(-4122)

$ws.Cells.Item(3, 2).Copy()
$ws.Range("H3:N3").PasteSpecial(-4122)

$ws.Cells.Item(4, 2).Copy()
$ws.Range("H4:N4").PasteSpecial(-4122)

$ws.Cells.Item(5, 2).Copy()
$ws.Range("H5:N5").PasteSpecial(-4122)

$excel.CutCopyMode = 0

# ---- Header row (row 1) ----------------------------------------------
$ws.Cells.Item(1, 2).Value  = "name"
$ws.Cells.Item(1, 3).Value  = "capacity"
$ws.Cells.Item(1, 4).Value  = "owner"
$ws.Cells.Item(1, 5).Value  = "register_date"
$ws.Cells.Item(1, 6).Value  = "register_reason"
$ws.Cells.Item(1, 7).Value  = "acquire_value"
$ws.Cells.Item(1, 8).Value  = "property_category"
$ws.Cells.Item(1, 9).Value  = "category"
$ws.Cells.Item(1, 10).Value = "date"
$ws.Cells.Item(1, 11).Value = "legislator_name"
$ws.Cells.Item(1, 12).Value = "legislator_id"
$ws.Cells.Item(1, 13).Value = "source_file"
$ws.Cells.Item(1, 14).Value = "index"

# ---- Metadata shared by every data row --------------------------------
$propertyCategory = "land"
$category         = "normal"
$recordDate       = "2013-05-01"
$legislatorName   = "顏寬恒"
$legislatorId     = 1803
$sourceFile       = "tmpbf3f1"

# Keep the record-date column as literal text (not an auto-converted
# date serial) for every data row.
$ws.Range("J2:J5").NumberFormat = "@"

# ---- Row 2 (index 112) ------------------------------------------------
$ws.Cells.Item(2, 8).Value  = $propertyCategory
$ws.Cells.Item(2, 9).Value  = $category
$ws.Cells.Item(2, 10).Value = $recordDate
$ws.Cells.Item(2, 11).Value = $legislatorName
$ws.Cells.Item(2, 12).Value = $legislatorId
$ws.Cells.Item(2, 13).Value = $sourceFile
$ws.Cells.Item(2, 14).Value = 112

# ---- Row 3 (index 113) ------------------------------------------------
$ws.Cells.Item(3, 8).Value  = $propertyCategory
$ws.Cells.Item(3, 9).Value  = $category
$ws.Cells.Item(3, 10).Value = $recordDate
$ws.Cells.Item(3, 11).Value = $legislatorName
$ws.Cells.Item(3, 12).Value = $legislatorId
$ws.Cells.Item(3, 13).Value = $sourceFile
$ws.Cells.Item(3, 14).Value = 113

# ---- Row 4 (index 114) ------------------------------------------------
$ws.Cells.Item(4, 8).Value  = $propertyCategory
$ws.Cells.Item(4, 9).Value  = $category
$ws.Cells.Item(4, 10).Value = $recordDate
$ws.Cells.Item(4, 11).Value = $legislatorName
$ws.Cells.Item(4, 12).Value = $legislatorId
$ws.Cells.Item(4, 13).Value = $sourceFile
$ws.Cells.Item(4, 14).Value = 114

# ---- Row 5 (index 115) ------------------------------------------------
$ws.Cells.Item(5, 8).Value  = $propertyCategory
$ws.Cells.Item(5, 9).Value  = $category
$ws.Cells.Item(5, 10).Value = $recordDate
$ws.Cells.Item(5, 11).Value = $legislatorName
$ws.Cells.Item(5, 12).Value = $legislatorId
$ws.Cells.Item(5, 13).Value = $sourceFile
$ws.Cells.Item(5, 14).Value = 115
